# =====================================================================
# "grafica del exp 2.2" - restructure Exp2-Entr2 sheet with a second
# data table (Con seguridad / Sin seguridad) and update the chart to
# show both series with trendlines, a legend, and new axis titles.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Sheet "Exp2-Entr2" (sheet4) - rebuild the data area
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Exp2-Entr2")

# --- Block 1: "Nodo con seguridad" (A1:D14) ---
$ws4.Range("A1:D1").Merge()
$ws4.Range("A1").Value = "Nodo con seguridad"

$ws4.Range("A2").Value = "# threads"
$ws4.Range("B2").Value = "Media (ms)"
$ws4.Range("C2").Value = "Error (%)"
$ws4.Range("D2").Value = "Throughput (threads/seg)"

$threadVals = @(1200,1100,1000,900,800,700,600,500,400,300,200,100)
for ($i = 0; $i -lt $threadVals.Count; $i++) {
    $r = 3 + $i
    $ws4.Range("A$r").Value = $threadVals[$i]
}

# B3:D14 -> "N/A" placeholders, except B5:B14 which keep the original
# "Media (ms)" series values (that used to live in B3:B12 / C3:C12)
$naRows = 3..14
foreach ($r in $naRows) {
    $ws4.Range("B$r").Value = "N/A"
    $ws4.Range("C$r").Value = "N/A"
    $ws4.Range("D$r").Value = "N/A"
}
$oldMedia = @{5=56012; 6=46232; 7=42384; 8=45934; 9=41143; 10=47429; 11=47283; 12=42649; 13=35609; 14=16911}
foreach ($r in $oldMedia.Keys) {
    $ws4.Range("B$r").Value = $oldMedia[$r]
}

# --- Block 2: "Gráfico" (F1:H12) ---
$ws4.Range("F1:H1").Merge()
$ws4.Range("F1").Value = "Gráfico"

$ws4.Range("F2").Value = "# thread"
$ws4.Range("G2").Value = "Con seguridad"
$ws4.Range("H2").Value = "Sin seguridad"

$gridF = @(1000,900,800,700,600,500,400,300,200,100)
$gridG = @(56012,46232,42384,45934,41143,47429,47283,42649,35609,16911)
$gridH = @(955,942,841,696,243,166,2,3,2,3)
for ($i = 0; $i -lt $gridF.Count; $i++) {
    $r = 3 + $i
    $ws4.Range("F$r").Value = $gridF[$i]
    $ws4.Range("G$r").Value = $gridG[$i]
    $ws4.Range("H$r").Value = $gridH[$i]
}

# --- Block 3: "Nodo único" (A15:D28) ---
$ws4.Range("A15:D15").Merge()
$ws4.Range("A15").Value = "Nodo único"

$ws4.Range("A16").Value = "# threads"
$ws4.Range("B16").Value = "Media (ms)"
$ws4.Range("C16").Value = "Error (%)"
$ws4.Range("D16").Value = "Throughput (threads/seg)"

$block3 = @(
    @(1200,1169,0,305.49898167006103),
    @(1100,984,0,289.32140978432398),
    @(1000,955,0,286.61507595299503),
    @(900,942,0,284.270372710044),
    @(800,841,0,290.06526468455399),
    @(700,696,0,297.99914857386102),
    @(600,243,0,308.64197530864197),
    @(500,166,0,318.06615776081401),
    @(400,2,0,266.31158455392801),
    @(300,3,0,294.11764705882302),
    @(200,2,0,173.91304347825999),
    @(100,3,0,87.7963125548727)
)
for ($i = 0; $i -lt $block3.Count; $i++) {
    $r = 17 + $i
    $row = $block3[$i]
    $ws4.Range("A$r").Value = $row[0]
    $ws4.Range("B$r").Value = $row[1]
    $ws4.Range("C$r").Value = $row[2]
    $ws4.Range("D$r").Value = $row[3]
}

# --- Formatting: headers bold + centered + thin border box ---
$headerRanges = @("A1:D1","F1:H1","A15:D15")
foreach ($hr in $headerRanges) {
    $rng = $ws4.Range($hr)
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = -4108
    $rng.Borders.LineStyle = 1
}

$colHeaderRanges = @("A2:D2","F2:H2","A16:D16")
foreach ($hr in $colHeaderRanges) {
    $rng = $ws4.Range($hr)
    $rng.Font.Bold = $true
    $rng.Borders.LineStyle = 1
}

# Outer box borders around each data block (left/right/bottom of the
# whole block + thin separators), matching the look of a bordered table
$boxRanges = @("A1:D14","F1:H12","A15:D28")
foreach ($br in $boxRanges) {
    $rng = $ws4.Range($br)
    $rng.BorderAround(1) | Out-Null
}
$ws4.Range("A1:D1").BorderAround(1) | Out-Null
$ws4.Range("F1:H1").BorderAround(1) | Out-Null
$ws4.Range("A15:D15").BorderAround(1) | Out-Null
$ws4.Range("A2:D2").BorderAround(1) | Out-Null
$ws4.Range("F2:H2").BorderAround(1) | Out-Null
$ws4.Range("A16:D16").BorderAround(1) | Out-Null

# D-column (Throughput) keeps the 0.0 numeric format on the numeric rows
$ws4.Range("D17:D28").NumberFormat = "0.0"

# Column widths (approximate best-fit)
$ws4.Columns.Item(4).ColumnWidth = 22.5
$ws4.Columns.Item(7).ColumnWidth = 12.66
$ws4.Columns.Item(8).ColumnWidth = 12

# Sheet view: scroll/zoom/selection + this sheet becomes the active tab
$ws4.Activate()
$excel.ActiveWindow.Zoom = 162
$ws4.Range("G1").Select()
$ws4.Range("O22").Select()

# ---------------------------------------------------------------
# 2. Chart "Chart 1" on Exp2-Entr2 -> rename to "Chart 2" and wire
#    the two series ("Con seguridad" / "Sin seguridad")
# ---------------------------------------------------------------
$co = $ws4.ChartObjects().Item(1)
$co.Name = "Chart 2"
$chart = $co.Chart

# Remove the chart title entirely (autoTitleDeleted = 1)
$chart.HasTitle = $false

# --- Series 1: "Con seguridad" (existing series, repointed) ---
$ser1 = $chart.SeriesCollection().Item(1)
$ser1.Name = "Con seguridad"
$ser1.XValues = "='Exp2-Entr2'!`$F`$3:`$F`$12"
$ser1.Values = "='Exp2-Entr2'!`$G`$3:`$G`$12"

$tl1 = $ser1.Trendlines().Add()
$tl1.Type = -4132
$tl1.DisplayEquation = $false
$tl1.DisplayRSquared = $false

# --- Series 2: "Sin seguridad" (new series) ---
$ser2 = $chart.SeriesCollection().NewSeries()
$ser2.Name = "Sin seguridad"
$ser2.XValues = "='Exp2-Entr2'!`$F`$3:`$F`$12"
$ser2.Values = "='Exp2-Entr2'!`$H`$3:`$H`$12"
$ser2.ChartType = 74

$tl2 = $ser2.Trendlines().Add()
$tl2.Type = -4132
$tl2.DisplayEquation = $false
$tl2.DisplayRSquared = $false

# --- Axes ---
$xAxis = $chart.Axes(1, 1)
$xAxis.HasTitle = $true
$xAxis.AxisTitle.Text = "Numero de threads"
$xAxis.MajorTickMark = 2
$xAxis.HasMinorGridlines = $true

$yAxis = $chart.Axes(2, 1)
$yAxis.HasTitle = $true
$yAxis.AxisTitle.Text = "Media (ms)"
$yAxis.MajorTickMark = 2
$yAxis.HasMinorGridlines = $true
$yAxis.MinimumScale = 0

# --- Legend ---
$chart.HasLegend = $true
$chart.Legend.Position = -4152

Write-Host "Exp2-Entr2 rebuilt + chart updated"

# ---------------------------------------------------------------
# 3. Sheet view / selection bookkeeping for the other sheets
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Exp2-Entr1(1)")
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 14
$ws2.Range("G22").Select()

$ws3 = $wb.Worksheets.Item("Exp2-Entr1(2)")
$ws3.Activate()
$excel.ActiveWindow.ScrollRow = 24
$ws3.Range("A36:D49").Select()

# Re-activate Exp2-Entr2 as the final active tab
$ws4.Activate()

Write-Host "done"
